$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "10/31/2025"
$ws.Range("A60").ClearFormats()
$ws.Range("B60").Value = 0.1918410875229571
$ws.Range("C60").Value = 0.8081589124770429
